# Scheduled runner update: refresh currentAveragePrice / leve profit columns
# (H..N) across the market-board snapshot rows that changed on this pass.
# Each entry is one cell on one sheet; Value=$null means the cell is cleared
# (the source row no longer carries that column, matching upstream source).

$wb = $excel.ActiveWorkbook

$updates = @(
    @{Sheet="ALC"; Cell="H2"; Value=1020.4706},
    @{Sheet="ALC"; Cell="J2"; Value=3359.6},
    @{Sheet="ALC"; Cell="L2"; Value=3359.6},
    @{Sheet="ALC"; Cell="N2"; Value=-3585.6},
    @{Sheet="ALC"; Cell="H6"; Value=322},
    @{Sheet="ALC"; Cell="J6"; Value=400},
    @{Sheet="ALC"; Cell="L6"; Value=1200},
    @{Sheet="ALC"; Cell="N6"; Value=-1424},
    @{Sheet="ALC"; Cell="H17"; Value=3525.739},
    @{Sheet="ALC"; Cell="J17"; Value=3604.55},
    @{Sheet="ALC"; Cell="L17"; Value=10813.65},
    @{Sheet="ALC"; Cell="N17"; Value=-11149.65},
    @{Sheet="ALC"; Cell="H34"; Value=0},
    @{Sheet="ALC"; Cell="I34"; Value=0},
    @{Sheet="ALC"; Cell="K34"; Value=0},
    @{Sheet="ALC"; Cell="M34"; Value=$null},
    @{Sheet="ALC"; Cell="H36"; Value=0},
    @{Sheet="ALC"; Cell="I36"; Value=0},
    @{Sheet="ALC"; Cell="K36"; Value=0},
    @{Sheet="ALC"; Cell="M36"; Value=$null},
    @{Sheet="ALC"; Cell="H70"; Value=4321.885},
    @{Sheet="ALC"; Cell="J70"; Value=4321.885},
    @{Sheet="ALC"; Cell="L70"; Value=12965.655},
    @{Sheet="ALC"; Cell="N70"; Value=-13505.655},
    @{Sheet="ALC"; Cell="H73"; Value=4321.885},
    @{Sheet="ALC"; Cell="J73"; Value=4321.885},
    @{Sheet="ALC"; Cell="L73"; Value=12965.655},
    @{Sheet="ALC"; Cell="N73"; Value=-14837.655},
    @{Sheet="ALC"; Cell="H98"; Value=1886.3889},
    @{Sheet="ALC"; Cell="I98"; Value=1613.5714},
    @{Sheet="ALC"; Cell="J98"; Value=2841.25},
    @{Sheet="ALC"; Cell="K98"; Value=1613.5714},
    @{Sheet="ALC"; Cell="L98"; Value=2841.25},
    @{Sheet="ALC"; Cell="M98"; Value=-115.5714},
    @{Sheet="ALC"; Cell="N98"; Value=-5837.25},
    @{Sheet="ALC"; Cell="H111"; Value=1106.4286},
    @{Sheet="ALC"; Cell="I111"; Value=799},
    @{Sheet="ALC"; Cell="J111"; Value=1875},
    @{Sheet="ALC"; Cell="K111"; Value=2397},
    @{Sheet="ALC"; Cell="L111"; Value=5625},
    @{Sheet="ALC"; Cell="M111"; Value=670},
    @{Sheet="ALC"; Cell="N111"; Value=-11759},
    @{Sheet="ALC"; Cell="H113"; Value=4713.2},
    @{Sheet="ALC"; Cell="J113"; Value=4855.6665},
    @{Sheet="ALC"; Cell="L113"; Value=4855.6665},
    @{Sheet="ALC"; Cell="N113"; Value=-11363.6665},
    @{Sheet="ALC"; Cell="H116"; Value=4733},
    @{Sheet="ALC"; Cell="I116"; Value=3944},
    @{Sheet="ALC"; Cell="K116"; Value=3944},
    @{Sheet="ALC"; Cell="M116"; Value=-502},
    @{Sheet="ALC"; Cell="H118"; Value=906},
    @{Sheet="ALC"; Cell="I118"; Value=897.63635},
    @{Sheet="ALC"; Cell="J118"; Value=998},
    @{Sheet="ALC"; Cell="K118"; Value=2692.90905},
    @{Sheet="ALC"; Cell="L118"; Value=2994},
    @{Sheet="ALC"; Cell="M118"; Value=-1035.90905},
    @{Sheet="ALC"; Cell="N118"; Value=-6308},
    @{Sheet="ALC"; Cell="H122"; Value=1886.3889},
    @{Sheet="ALC"; Cell="I122"; Value=1613.5714},
    @{Sheet="ALC"; Cell="J122"; Value=2841.25},
    @{Sheet="ALC"; Cell="K122"; Value=4840.7142},
    @{Sheet="ALC"; Cell="L122"; Value=8523.75},
    @{Sheet="ALC"; Cell="M122"; Value=-2390.7142},
    @{Sheet="ALC"; Cell="N122"; Value=-13423.75},
    @{Sheet="ALC"; Cell="H132"; Value=5990},
    @{Sheet="ALC"; Cell="I132"; Value=0},
    @{Sheet="ALC"; Cell="J132"; Value=5990},
    @{Sheet="ALC"; Cell="K132"; Value=0},
    @{Sheet="ALC"; Cell="L132"; Value=17970},
    @{Sheet="ALC"; Cell="M132"; Value=$null},
    @{Sheet="ALC"; Cell="N132"; Value=-23030},
    @{Sheet="ARM"; Cell="H32"; Value=15771.9},
    @{Sheet="ARM"; Cell="I32"; Value=14191},
    @{Sheet="ARM"; Cell="K32"; Value=14191},
    @{Sheet="ARM"; Cell="M32"; Value=-13904},
    @{Sheet="BSM"; Cell="H20"; Value=2560.375},
    @{Sheet="BSM"; Cell="I20"; Value=2003.5},
    @{Sheet="BSM"; Cell="K20"; Value=2003.5},
    @{Sheet="BSM"; Cell="M20"; Value=-1756.5},
    @{Sheet="BSM"; Cell="H86"; Value=3692.7},
    @{Sheet="BSM"; Cell="I86"; Value=2329},
    @{Sheet="BSM"; Cell="J86"; Value=5738.25},
    @{Sheet="BSM"; Cell="K86"; Value=2329},
    @{Sheet="BSM"; Cell="L86"; Value=5738.25},
    @{Sheet="BSM"; Cell="M86"; Value=-1206},
    @{Sheet="BSM"; Cell="N86"; Value=-7984.25},
    @{Sheet="BSM"; Cell="H89"; Value=3692.7},
    @{Sheet="BSM"; Cell="I89"; Value=2329},
    @{Sheet="BSM"; Cell="J89"; Value=5738.25},
    @{Sheet="BSM"; Cell="K89"; Value=11645},
    @{Sheet="BSM"; Cell="L89"; Value=28691.25},
    @{Sheet="BSM"; Cell="M89"; Value=-6029},
    @{Sheet="BSM"; Cell="N89"; Value=-39923.25},
    @{Sheet="BSM"; Cell="H94"; Value=1201.0834},
    @{Sheet="BSM"; Cell="I94"; Value=1327.0625},
    @{Sheet="BSM"; Cell="J94"; Value=949.125},
    @{Sheet="BSM"; Cell="K94"; Value=1327.0625},
    @{Sheet="BSM"; Cell="L94"; Value=949.125},
    @{Sheet="BSM"; Cell="M94"; Value=-876.0625},
    @{Sheet="BSM"; Cell="N94"; Value=-1851.125},
    @{Sheet="BSM"; Cell="H117"; Value=0},
    @{Sheet="BSM"; Cell="J117"; Value=0},
    @{Sheet="BSM"; Cell="L117"; Value=0},
    @{Sheet="BSM"; Cell="N117"; Value=$null},
    @{Sheet="CUL"; Cell="H92"; Value=500},
    @{Sheet="CUL"; Cell="I92"; Value=500},
    @{Sheet="CUL"; Cell="K92"; Value=1500},
    @{Sheet="CUL"; Cell="M92"; Value=-252},
    @{Sheet="CUL"; Cell="H107"; Value=506.9},
    @{Sheet="CUL"; Cell="I107"; Value=274.08334},
    @{Sheet="CUL"; Cell="K107"; Value=822.2500200000001},
    @{Sheet="CUL"; Cell="M107"; Value=1097.74998},
    @{Sheet="CUL"; Cell="H123"; Value=3000},
    @{Sheet="CUL"; Cell="I123"; Value=3000},
    @{Sheet="CUL"; Cell="K123"; Value=9000},
    @{Sheet="CUL"; Cell="M123"; Value=-6550},
    @{Sheet="GSM"; Cell="H2"; Value=469},
    @{Sheet="GSM"; Cell="I2"; Value=459.33334},
    @{Sheet="GSM"; Cell="J2"; Value=498},
    @{Sheet="GSM"; Cell="K2"; Value=459.33334},
    @{Sheet="GSM"; Cell="L2"; Value=498},
    @{Sheet="GSM"; Cell="M2"; Value=-346.33334},
    @{Sheet="GSM"; Cell="N2"; Value=-724},
    @{Sheet="GSM"; Cell="H97"; Value=1616.2},
    @{Sheet="GSM"; Cell="I97"; Value=1616.2},
    @{Sheet="GSM"; Cell="J97"; Value=0},
    @{Sheet="GSM"; Cell="K97"; Value=1616.2},
    @{Sheet="GSM"; Cell="L97"; Value=0},
    @{Sheet="GSM"; Cell="M97"; Value=-1120.2},
    @{Sheet="GSM"; Cell="N97"; Value=$null},
    @{Sheet="LTW"; Cell="H16"; Value=12998.5},
    @{Sheet="LTW"; Cell="I16"; Value=8998},
    @{Sheet="LTW"; Cell="K16"; Value=8998},
    @{Sheet="LTW"; Cell="M16"; Value=-8828},
    @{Sheet="LTW"; Cell="H40"; Value=7419.3335},
    @{Sheet="LTW"; Cell="I40"; Value=7396.2856},
    @{Sheet="LTW"; Cell="K40"; Value=7396.2856},
    @{Sheet="LTW"; Cell="M40"; Value=-7260.2856},
    @{Sheet="LTW"; Cell="H82"; Value=1557.4},
    @{Sheet="LTW"; Cell="I82"; Value=893.75},
    @{Sheet="LTW"; Cell="K82"; Value=893.75},
    @{Sheet="LTW"; Cell="M82"; Value=-532.75},
    @{Sheet="LTW"; Cell="H85"; Value=1557.4},
    @{Sheet="LTW"; Cell="I85"; Value=893.75},
    @{Sheet="LTW"; Cell="K85"; Value=893.75},
    @{Sheet="LTW"; Cell="M85"; Value=354.25},
    @{Sheet="LTW"; Cell="H136"; Value=3169.8333},
    @{Sheet="LTW"; Cell="I136"; Value=3169.8333},
    @{Sheet="LTW"; Cell="K136"; Value=9509.499899999999},
    @{Sheet="LTW"; Cell="M136"; Value=-6959.499899999999},
    @{Sheet="WVR"; Cell="H62"; Value=2421.75},
    @{Sheet="WVR"; Cell="I62"; Value=2462.3333},
    @{Sheet="WVR"; Cell="J62"; Value=2300},
    @{Sheet="WVR"; Cell="K62"; Value=2462.3333},
    @{Sheet="WVR"; Cell="L62"; Value=2300},
    @{Sheet="WVR"; Cell="M62"; Value=-1838.3333},
    @{Sheet="WVR"; Cell="N62"; Value=-3548},
    @{Sheet="WVR"; Cell="H65"; Value=2421.75},
    @{Sheet="WVR"; Cell="I65"; Value=2462.3333},
    @{Sheet="WVR"; Cell="J65"; Value=2300},
    @{Sheet="WVR"; Cell="K65"; Value=12311.6665},
    @{Sheet="WVR"; Cell="L65"; Value=11500},
    @{Sheet="WVR"; Cell="M65"; Value=-9191.666499999999},
    @{Sheet="WVR"; Cell="N65"; Value=-17740},
    @{Sheet="WVR"; Cell="H136"; Value=1681.5},
    @{Sheet="WVR"; Cell="I136"; Value=1594.9048},
    @{Sheet="WVR"; Cell="K136"; Value=4784.7144},
    @{Sheet="WVR"; Cell="M136"; Value=-2234.7144}
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $cell = $ws.Range($u.Cell)
    if ($null -eq $u.Value) {
        $cell.ClearContents()
    } else {
        $cell.Value = $u.Value
    }
}
